# Applies the "Added files via upload" commit:
#   - Inserts two new slides (Title and Content layout) at positions 4 and 5:
#       4: "Division of Work/Community Formation/Sustainability"
#       5: "Schedule"
#     (the former slides 4-7, all full-bleed screenshots, shift down to 6-9)

$p = $ppt.ActivePresentation

$rsq  = [char]0x2019   # RIGHT SINGLE QUOTATION MARK  ( ' )
$endash = [char]0x2013 # EN DASH ( - )
$tab  = [char]9

$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"

# ---------------------------------------------------------------------
# Insert in this order so the internally generated slide IDs land as the
# source deck has them: id 263 -> "Schedule" (position 5), then
# id 264 -> "Division of Work..." (position 4, pushing Schedule to 5).
# ---------------------------------------------------------------------

# --- New slide at position 4 (will become "Schedule" first, then get
#     pushed to slot 5 once the Division slide is inserted before it) ---
$sSchedule = $p.Slides.AddSlide(4, $titleContentLayout)

$scheduleBody = $sSchedule.Shapes.Placeholders.Item(1)
$tr = $scheduleBody.TextFrame.TextRange
$tr.Text = "Weeks 1 " + $endash + " 3 : Learn "
$tr.InsertAfter("Ocaml") | Out-Null
$tr.InsertAfter(" and get familiar with the environment") | Out-Null
$tr.InsertAfter("`rWeeks 4-6: Start adding the project and working on the goals we set out to achieve. Add items/") | Out-Null
$tr.InsertAfter("npcs") | Out-Null
$tr.InsertAfter("/moves to the game") | Out-Null
$tr.InsertAfter("`rWeek 7: Review our code and get ready to submit back to the main repository") | Out-Null

$scheduleTitle = $sSchedule.Shapes.Placeholders.Item(2)
$scheduleTitle.TextFrame.TextRange.Text = "Schedule" + $tab

# --- New slide at position 4: "Division of Work/Community Formation/Sustainability" ---
$sDivision = $p.Slides.AddSlide(4, $titleContentLayout)

$divisionBody = $sDivision.Shapes.Placeholders.Item(1)
$tr = $divisionBody.TextFrame.TextRange
$tr.Text = "We will separately work on each category we are adding, "
$tr.InsertAfter("ie") | Out-Null
$tr.InsertAfter(" moves/items/NPC" + $rsq + "s") | Out-Null
$tr.InsertAfter("`rWe expect to be able to work with the project" + $rsq + "s creators if any problems arise") | Out-Null
$tr.InsertAfter("`rAt the end we will send a pull request to the master branch") | Out-Null
$tr.InsertAfter("`rThis project will be sustainable as we are filling in the missing information and there is always more features to add") | Out-Null
$tr.InsertAfter("`rThey are trying to port this game to windows so our changes will be transferred to that ") | Out-Null
$tr.InsertAfter("new version") | Out-Null
$divisionBody.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>

$divisionTitle = $sDivision.Shapes.Placeholders.Item(2)
$divisionTitle.TextFrame.TextRange.Text = "Division of Work/Community Formation/Sustainability" + $tab + $tab
$divisionTitle.TextFrame.AutoSize = 2   # ppAutoSizeTextToFitShape -> <a:normAutofit/>
